# Calcular Multas (Com variável Multas como Aleatória)
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) Dados_Projetados: header M1 becomes a formula; K/L/M data cols updated
# ---------------------------------------------------------------
$wsDados = $wb.Worksheets.Item("Dados_Projetados")

$wsDados.Range("M1").Formula = '="CustoMedioMulta_Lei1"'

$wsDados.Range("K2").Value = 1
$wsDados.Range("L2").Value = 0
$wsDados.Range("M2").Value = 1200

$wsDados.Range("K3").Value = 1
$wsDados.Range("L3").Value = 0
$wsDados.Range("M3").Value = 1200

$wsDados.Range("L4").Value = 0
$wsDados.Range("M4").Value = 1200

$wsDados.Range("L5").Value = 0
$wsDados.Range("M5").Value = 1200

$wsDados.Range("L6").Value = 0
$wsDados.Range("M6").Value = 1200

$wsDados.Range("L3:L6").Select()

# ---------------------------------------------------------------
# 2) Parametros: rows 24-27 renamed / retargeted, rows 28-41 removed
# ---------------------------------------------------------------
$wsParam = $wb.Worksheets.Item("Parametros")

$wsParam.Range("A24").Value = "NumeroMultasAPriori_Lei1"
$wsParam.Range("C24").Value = 2
$wsParam.Range("D24").Value = 0.5

$wsParam.Range("A25").Value = "Atendimento_Lei1"
$wsParam.Range("C25").Value = 0
$wsParam.Range("D25").Value = 0

$wsParam.Range("A26").Value = "NumeroMultasAPriori_Lei1"
$wsParam.Range("C26").Value = 2
$wsParam.Range("D26").Value = 0.5
$wsParam.Range("G26").Value = "Iniciativa1"

$wsParam.Range("A27").Value = "Atendimento_Lei1"
$wsParam.Range("C27").Value = 1
$wsParam.Range("D27").Value = 0
$wsParam.Range("G27").Value = "Iniciativa1"

$wsParam.Range("A28:A41").EntireRow.Delete()

# ---------------------------------------------------------------
# 3) Funcoes_Inputs: rows 31-35 retargeted, rows 36-50 removed
# ---------------------------------------------------------------
$wsFuncIn = $wb.Worksheets.Item("Funcoes_Inputs")

$wsFuncIn.Range("B31").Formula = '="NumeroMultasAPriori_"&D31'

$wsFuncIn.Range("B32").Formula = '="CustoMedioMulta_"&D32'
$wsFuncIn.Range("C32").Formula = "=FALSE"
$wsFuncIn.Range("D32").Value = "Lei1"

$wsFuncIn.Range("B33").Formula = '="Atendimento_"&D33'
$wsFuncIn.Range("D33").Value = "Lei1"

$wsFuncIn.Range("B34").Value = "Crise"
$wsFuncIn.Range("C34").Formula = "=FALSE"
$wsFuncIn.Range("D34").ClearContents()

$wsFuncIn.Range("B35").Value = "FatorCrise"
$wsFuncIn.Range("C35").Formula = "=FALSE"
$wsFuncIn.Range("D35").ClearContents()

$wsFuncIn.Range("A36:A50").EntireRow.Delete()

$wsFuncIn.Range("B31,B33").Select()
$wsFuncIn.Range("B33").Activate()

# ---------------------------------------------------------------
# 4) Funcoes_Outputs: rows 25-29 removed
# ---------------------------------------------------------------
$wsFuncOut = $wb.Worksheets.Item("Funcoes_Outputs")

$wsFuncOut.Range("A25:A29").EntireRow.Delete()

$wsFuncOut.Range("B26").Select()

# ---------------------------------------------------------------
# Restore Parametros as the active sheet/selection (matches workbook.xml
# bookViews activeTab, which the commit does not change)
# ---------------------------------------------------------------
$wsParam.Select()
$wsParam.Range("F20").Select()
